# Added BAARD stage 2, line, watermark to CIFAR10
# The "apgd2 @ 0.2" data row (row 8) is removed; all subsequent rows shift up
# by one, and the merged "apgd2" label moves up with the rest of its group.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(8).Delete()

# Re-apply the group label that was attached to the now-deleted physical row;
# it belongs to the merged range that now starts at A8.
$ws.Range("A8").Value = "apgd2"
